$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting of existing data rows (2-9) down through the new rows (10-24)
# so new rows inherit the same cell styles (e.g. the timestamp number format on column D).
$ws.Range("A2:D9").Copy() | Out-Null
$ws.Range("A10:D17").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:D8").Copy() | Out-Null
$ws.Range("A18:D24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$data = @(
    ,@(2, '0ce5dd49', 'Summarising lecture notes or readings|Generating practice questions or quizzes|Explaining difficult concepts in simple terms|Reviewing flashcards – key terms', 'q11_study_subtasks', 45854.6558148386)
    ,@(3, '2c1001cb', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(4, '37cc37bf', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(5, '43faa0b9', 'Summarising lecture notes or readings|Generating practice questions or quizzes|Explaining difficult concepts in simple terms|Creating mnemonics or memory aids|Reviewing flashcards – key terms', 'q11_study_subtasks', 45854.6558148386)
    ,@(6, '4abe3e88', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(7, '50164f59', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(8, '5cf70f79', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(9, '5da96769', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(10, '6ca3e2f6', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(11, '790a4fcb', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(12, '802cc63a', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(13, '85c3ea4d', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(14, '942dfafb', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(15, '9bc6ba8c', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(16, 'a2d65af2', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(17, 'a46f1771', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(18, 'ad58f9da', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(19, 'c7d9a301', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(20, 'ce8732ff', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(21, 'd6f1d567', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(22, 'da9326c9', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(23, 'e09ca7bf', 'I did not choose “Study revision / exam prep”', 'q11_study_subtasks', 45854.6558148386)
    ,@(24, 'ef53a641', 'Summarising lecture notes or readings|Generating practice questions or quizzes|Explaining difficult concepts in simple terms|Reviewing flashcards – key terms', 'q11_study_subtasks', 45854.6558148386)
)

foreach ($rec in $data) {
    $r = $rec[0]
    $ws.Cells.Item($r, 1).Value2 = $rec[1]
    $ws.Cells.Item($r, 2).Value2 = $rec[2]
    $ws.Cells.Item($r, 3).Value2 = $rec[3]
    $ws.Cells.Item($r, 4).Value2 = [double]$rec[4]
}

